# adding averages and more checks
#
# 1) Header / title styling: the bold-size-14 "title" font and the bold
#    "header row" font are consolidated into a single bold + white font
#    (used by both the dashboard title cell and the column-header row on
#    both sheets).
# 2) Training Dashboard ("PERIOD TO EXPIRE" / "LAST UPDATE" columns H & I,
#    rows 3-32): the report was regenerated 8 days later (LAST UPDATE moves
#    from 08-Sep-2025 to 16-Sep-2025), so the remaining days-to-expiry in
#    column H drops by 8 for every row.
# 3) Exam Dashboard: the "COMMENTS" column (E3:E10) wording changes from
#    "OK" to "date is valid", and its column gets a bit wider to fit the
#    new text.

$wb = $excel.ActiveWorkbook

$trainingWs = $wb.Worksheets.Item("Training Dashboard")
$examWs = $wb.Worksheets.Item("Exam Dashboard")

# --- 1) Title + header-row font: bold, white text (no more oversized title
# font) on both sheets --------------------------------------------------
foreach ($ws in @($trainingWs, $examWs)) {
    $title = $ws.Range("A1")
    $title.Font.Bold = $true
    $title.Font.Size = 11
    $title.Font.Color = 0xFFFFFF

    $headerRow = $ws.Rows.Item(2)
    $headerRow.Font.Bold = $true
    $headerRow.Font.Color = 0xFFFFFF
}

# --- 2) Training Dashboard: refresh PERIOD TO EXPIRE / LAST UPDATE -----
for ($r = 3; $r -le 32; $r++) {
    $periodCell = $trainingWs.Cells.Item($r, 8)
    $periodCell.Value = $periodCell.Value2 - 8
    # leading apostrophe forces text (matches the sheet's literal-text
    # storage instead of Excel auto-converting the string to a date serial)
    $trainingWs.Cells.Item($r, 9).Value = "'16-Sep-2025"
}

# --- 3) Exam Dashboard: comments column wording + column width ---------
$examWs.Columns.Item(5).ColumnWidth = 14.17

for ($r = 3; $r -le 10; $r++) {
    $examWs.Cells.Item($r, 5).Value = "date is valid"
}
